$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.648.19"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "1.631.45"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'212.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'0.253"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.70%  "

$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("D10").Value = "'19.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.83%  "

$ws.Range("D11").Value = "'0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.21%  "

$ws.Range("D12").Value = "1.858.50"
$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").Value = "1.628.51"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D16").Value = "26.629.04"
$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("D17").Value = "'63.32"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "0.0₃0743"
$ws.Range("E18").Value = "  +2.14%  "

$ws.Range("D19").Value = "'218.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.66%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "'4.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").Value = "'6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.27%  "

$ws.Range("D23").Value = "'9.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("E24").Value = "  +4.20%  "

$ws.Range("D25").Value = "'147.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("E28").Value = "  +3.97%  "

$ws.Range("D29").Value = "'15.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.40%  "

$ws.Range("D30").Value = "'0.0503"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.59%  "

$ws.Range("E31").Value = "  +0.35%  "

$ws.Range("E32").Value = "  +3.86%  "

$ws.Range("E33").Value = "  +2.36%  "

$ws.Range("E34").Value = "  +0.68%  "

$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "1.214.06"
$ws.Range("E36").Value = "  +3.04%  "

$ws.Range("D37").Value = "'0.0171"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.39%  "

$ws.Range("D38").Value = "'0.807"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("E41").Value = "  -2.11%  "

$ws.Range("D43").Value = "'0.792"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.82%  "

$ws.Range("D44").Value = "1.767.73"
$ws.Range("E44").Value = "  +0.43%  "

$ws.Range("D45").Value = "'92.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("E46").Value = "  +1.86%  "

$ws.Range("D47").Value = "'54.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.95%  "

$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("E49").Value = "  +3.87%  "

$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("E51").Value = "  +0.18%  "
